$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window view sizing (workbook-level bookViews) ---
$excel.Width = 20490
$excel.Height = 7620

# --- Sheet view: scrolled position + active selection ---
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D16").Select()

# --- Cell value updates (Ideal/Actual Tasks Remaining) ---
$ws.Range("B13").Value = 36

$ws.Range("C13").Value = 36
$ws.Range("C14").Value = 33
$ws.Range("C15").Value = 30
$ws.Range("C16").Value = 26
$ws.Range("C17").Value = 22
$ws.Range("C18").Value = 21
$ws.Range("C19").Value = 20
$ws.Range("C20").Value = 19
$ws.Range("C21").Value = 18
$ws.Range("C22").Value = 18
$ws.Range("C23").Value = 18
$ws.Range("C24").Value = 16
$ws.Range("C25").Value = 15
$ws.Range("C26").Value = 11
$ws.Range("C27").Value = 11
$ws.Range("C28").Value = 10
$ws.Range("C29").Value = 9
$ws.Range("C30").Value = 5
$ws.Range("C31").Value = 2
$ws.Range("C32").Value = 0
